# Pop the first unused name off "Sheet1" and record it as newly-used on
# the "used" log sheet.

$wb = $excel.ActiveWorkbook
$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# Grab the first available name (row 1, column A) before removing it.
# Use Value2 (not Value) to avoid a COM-variant quirk when the property
# is read through an intermediate variable.
$nextId = $namesSheet.Range("A1").Value2

# Remove it from the pool; Excel shifts the remaining rows up by one.
$namesSheet.Rows.Item(1).Delete()

# Append a new row to the "used" log describing when/where it was used.
$usedRow = $usedSheet.Range("A1").End(-4121).Row + 1

$usedSheet.Cells.Item($usedRow, 1).Value = $nextId
$usedSheet.Cells.Item($usedRow, 2).Value = "ChatGPT Image 2026年1月18日 11_38_38.png"
$usedSheet.Cells.Item($usedRow, 3).Value = "2026-01-18 11:40:05"
